# health_life_expectancy.xlsx -- "Text updates as supplied by PM&C."
#
# 1. Tweak wording in the Description sheet's narrative paragraph (B6):
#    "see our chapter on" -> "see the page on"
# 2. Append a new "Source" row (row 9) to the Description sheet with the
#    ABS citation.
# 3. Leave the Description sheet selected/active (it becomes the active
#    tab instead of Data).

$wb = $excel.ActiveWorkbook

$wsDesc = $wb.Worksheets.Item("Description")

# --- 1. Wording tweak on the Description sheet -----------------------------
$wsDesc.Range("B6").Value = "All states and territories have shown a small increase over the same period, with the exception of females in the Northern Territory. Indigenous life expectancy continues to show a substantial gap compared to non-Indigenous Australians (see the page on the National Indigenous Reform Agreement for more information). "

# --- 2. New Source row -------------------------------------------------
$wsDesc.Range("A9").Value = "Source"
$wsDesc.Range("B9").Value = "ABS 2016 and previous years, Life Tables, Australia, States and Territories, various years."

# Match the look of the other wrapped-paragraph rows (B5/B6/B7): 12pt black
# Arial, wrapped -- and give the row the same kind of explicit height those
# rows carry.
$wsDesc.Range("B9").Font.Name = "Arial"
$wsDesc.Range("B9").Font.Size = 12
$wsDesc.Range("B9").Font.Color = 0
$wsDesc.Range("B9").WrapText = $true
$wsDesc.Rows.Item(9).RowHeight = 26.95

# --- 3. Selection / active sheet ----------------------------------------
# Data's own selection (B12) is untouched; Description becomes the active
# sheet with its new Source row selected.
$wsDesc.Range("B9").Select()
$wsDesc.Activate()
